$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
"col10 width: " + $ws1.Columns.Item(10).ColumnWidth
"col1 width: " + $ws1.Columns.Item(1).ColumnWidth
"col3 width: " + $ws1.Columns.Item(3).ColumnWidth
